# feat: add 2022-Q1 data
#
# Insert a new "2022-Q1" sheet (fund-holdings detail) right before the
# existing "总计" (totals) summary sheet, and prepend a corresponding
# 2022-Q1 row to the "总计" sheet's history table.

$wb = $excel.ActiveWorkbook

# --- create the new "2022-Q1" sheet just before "总计" --------------------
# NOTE: worksheet variables captured before a structural change (sheet
# Add/Delete/Move) track POSITION, not identity, so any sheet reference we
# plan to use later must be re-fetched (by name) only after all Add calls
# have completed.
$q1 = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$q1.Name = "2022-Q1"

# Re-fetch the totals sheet now that the sheet collection is stable again.
$totalSheet = $wb.Worksheets.Item("总计")

# --- "2022-Q1" sheet: header row (row 1, columns B:H) ----------------------
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Match the bold/boxed header style + centred index-column style used by
# every other quarter sheet (e.g. "2021-Q4") by copying their formatting.
$styleSrc = $wb.Worksheets.Item("2021-Q4")
$styleSrc.Range("B1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$styleSrc.Range("A2").Copy()
$q1.Range("A2:A7").PasteSpecial(-4122)

# Fund holdings data (rows 2-7). Column A is a zero-based row index;
# columns B and D-G hold values that must stay TEXT (fund codes with
# leading zeros, and numeric-looking percentages/amounts stored as text,
# matching the other quarter sheets) - entered with a leading apostrophe to
# force text rather than numeric interpretation. Column H (仓位排名) is a
# genuine number.
$rows = @(
    @(0, "213001", "宝盈鸿利收益灵活配置混合A",             "17.98", "90.37", "3.78", "0.6796", 9),
    @(1, "011756", "博时产业优选灵活配置混合型证券投资基金A", "28.12", "67.04", "2.24", "0.6299", 7),
    @(2, "007581", "宝盈鸿利收益灵活配置混合C",             "0.73",  "90.37", "3.78", "0.0276", 9),
    @(3, "005482", "博时创新驱动灵活配置混合A",             "0.41",  "88.65", "4.96", "0.0203", 5),
    @(4, "011757", "博时产业优选灵活配置混合型证券投资基金C", "0.84",  "67.04", "2.24", "0.0188", 7),
    @(5, "005483", "博时创新驱动灵活配置混合C",             "0.05",  "88.65", "4.96", "0.0025", 5)
)

$r = 2
foreach ($row in $rows) {
    $q1.Cells.Item($r, 1).Value = $row[0]
    $q1.Cells.Item($r, 2).Value = "'" + $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $q1.Cells.Item($r, 4).Value = "'" + $row[3]
    $q1.Cells.Item($r, 5).Value = "'" + $row[4]
    $q1.Cells.Item($r, 6).Value = "'" + $row[5]
    $q1.Cells.Item($r, 7).Value = "'" + $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# --- prepend the 2022-Q1 summary row to "总计" -----------------------------
# Shift existing data rows (2-6) down by one row (to 3-7) first, preserving
# their date/count/value, then write the new 2022-Q1 row into row 2.
# NOTE: read via .Value2 (.Value's getter is not reliable in this runtime -
# it can yield a stub description string instead of the cell's contents).
for ($i = 6; $i -ge 2; $i--) {
    $totalSheet.Cells.Item($i + 1, 1).Value = $i - 1
    $totalSheet.Cells.Item($i + 1, 2).Value = $totalSheet.Cells.Item($i, 2).Value2
    $totalSheet.Cells.Item($i + 1, 3).Value = $totalSheet.Cells.Item($i, 3).Value2
    $totalSheet.Cells.Item($i + 1, 4).Value = $totalSheet.Cells.Item($i, 4).Value2
}

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 6
$totalSheet.Cells.Item(2, 4).Value = 1.38
